$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.400.83'
$ws.Range('E2').Value = '  -2.42%  '
$ws.Range('D3').Value = '2.960.89'
$ws.Range('E3').Value = '  -1.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.96'
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.14'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '2.954.60'
$ws.Range('E8').Value = '  -1.83%  '
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.10'
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.432'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000216'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.72'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '3.447.01'
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '60.642.26'
$ws.Range('E17').Value = '  -1.97%  '
$ws.Range('D18').Value = '2.968.78'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.41'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '451.83'
$ws.Range('E20').Value = '  -3.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.87'
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.661'
$ws.Range('E22').Value = '  -2.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.72'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.29'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.62'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.61'
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('B30').Value = 'Mantle'
$ws.Range('C30').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.12'
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '24.86'
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.81'
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '55.09'
$ws.Range('E33').Value = '  -2.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.23'
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.26'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.67'
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '446.62'
$ws.Range('E37').Value = '  -3.33%  '
$ws.Range('D38').Value = '3.147.67'
$ws.Range('E38').Value = '  +3.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0767'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0373'
$ws.Range('E40').Value = '  -2.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.115'
$ws.Range('E41').Value = '  +3.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.90'
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.38'
$ws.Range('E43').Value = '  -3.95%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.240'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.98'
$ws.Range('E46').Value = '  +5.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '117.31'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.106'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.91'
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('D50').Value = '0.0₃0499'
$ws.Range('E50').Value = '  -4.13%  '
$ws.Range('E51').Value = '  +8.35%  '
